$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while keeping the cell as TEXT, matching the
# source feed's formatting (e.g. trailing zeros like "0.900", "1.00").
# Excel auto-detects numeric-looking strings when assigned via .Value,
# so we briefly force a text number format, assign, then clear the
# format back to the sheet default (no numFmt is left applied to the
# cell afterwards).
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

# Row-by-row Price (D) / Volume(1h) (E) refresh
Set-TextValue $ws.Range("D2") "37.060.39"
Set-TextValue $ws.Range("E2") "  +0.48%  "
Set-TextValue $ws.Range("D3") "2.042.81"
Set-TextValue $ws.Range("E3") "  +0.09%  "
Set-TextValue $ws.Range("E4") "  -0.21%  "
Set-TextValue $ws.Range("D5") "247.34"
Set-TextValue $ws.Range("E5") "  -1.37%  "
Set-TextValue $ws.Range("D6") "0.663"
Set-TextValue $ws.Range("E6") "  -0.32%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "56.56"
Set-TextValue $ws.Range("E8") "  -1.70%  "
Set-TextValue $ws.Range("D9") "0.383"
Set-TextValue $ws.Range("E9") "  -0.52%  "
Set-TextValue $ws.Range("D10") "0.0779"
Set-TextValue $ws.Range("E10") "  -0.41%  "
Set-TextValue $ws.Range("E11") "  +0.42%  "
Set-TextValue $ws.Range("D12") "15.98"
Set-TextValue $ws.Range("E12") "  -2.02%  "
Set-TextValue $ws.Range("D13") "0.900"
Set-TextValue $ws.Range("E13") "  +12.24%  "
Set-TextValue $ws.Range("D14") "2.342.57"
Set-TextValue $ws.Range("E14") "  +0.12%  "
Set-TextValue $ws.Range("E15") "  +2.85%  "
Set-TextValue $ws.Range("D16") "2.043.60"
Set-TextValue $ws.Range("E16") "  +0.02%  "
Set-TextValue $ws.Range("D17") "19.05"
Set-TextValue $ws.Range("E17") "  +14.55%  "
Set-TextValue $ws.Range("D18") "37.054.42"
Set-TextValue $ws.Range("E18") "  +0.60%  "
Set-TextValue $ws.Range("D19") "74.76"
Set-TextValue $ws.Range("E19") "  -0.58%  "
Set-TextValue $ws.Range("D20") "0.0₃0890"
Set-TextValue $ws.Range("E20") "  -1.61%  "
Set-TextValue $ws.Range("D21") "5.40"
Set-TextValue $ws.Range("E21") "  +0.21%  "
Set-TextValue $ws.Range("D22") "236.65"
Set-TextValue $ws.Range("E22") "  -0.08%  "
Set-TextValue $ws.Range("D23") "1.00"
Set-TextValue $ws.Range("E24") "  +4.09%  "
Set-TextValue $ws.Range("D25") "9.56"
Set-TextValue $ws.Range("E25") "  +3.87%  "
Set-TextValue $ws.Range("D26") "171.05"
Set-TextValue $ws.Range("E26") "  +1.31%  "
Set-TextValue $ws.Range("D27") "2.17"
Set-TextValue $ws.Range("E27") "  -7.01%  "
Set-TextValue $ws.Range("D28") "20.10"
Set-TextValue $ws.Range("E28") "  -0.15%  "
Set-TextValue $ws.Range("E29") "  -0.65%  "
Set-TextValue $ws.Range("D32") "0.0620"
Set-TextValue $ws.Range("E32") "  +0.32%  "
Set-TextValue $ws.Range("D33") "4.62"
Set-TextValue $ws.Range("E33") "  +4.35%  "
Set-TextValue $ws.Range("D34") "0.0877"
Set-TextValue $ws.Range("E34") "  -0.64%  "
Set-TextValue $ws.Range("E35") "  -0.05%  "
Set-TextValue $ws.Range("E36") "  +4.99%  "
Set-TextValue $ws.Range("E37") "  +2.53%  "
Set-TextValue $ws.Range("E38") "  -1.54%  "
Set-TextValue $ws.Range("D41") "0.0995"
Set-TextValue $ws.Range("E41") "  -10.97%  "
Set-TextValue $ws.Range("D42") "0.0222"
Set-TextValue $ws.Range("E42") "  -0.20%  "
Set-TextValue $ws.Range("E43") "  +1.77%  "
Set-TextValue $ws.Range("D44") "17.14"
Set-TextValue $ws.Range("E44") "  -3.04%  "
Set-TextValue $ws.Range("D45") "97.10"
Set-TextValue $ws.Range("E45") "  +0.67%  "
Set-TextValue $ws.Range("D46") "2.39"
Set-TextValue $ws.Range("E46") "  -3.48%  "
Set-TextValue $ws.Range("D47") "1.279.45"
Set-TextValue $ws.Range("E47") "  +0.05%  "
Set-TextValue $ws.Range("D48") "2.84"
Set-TextValue $ws.Range("E48") "  -1.56%  "
Set-TextValue $ws.Range("D49") "6.78"
Set-TextValue $ws.Range("E49") "  +0.98%  "
Set-TextValue $ws.Range("D50") "2.227.70"
Set-TextValue $ws.Range("E50") "  -0.25%  "
Set-TextValue $ws.Range("D51") "44.46"
Set-TextValue $ws.Range("E51") "  +1.73%  "

# Row 30 <-> Row 31: Filecoin and ImmutableX swapped ranking order, with
# refreshed price/volume figures.
Set-TextValue $ws.Range("B30") "Filecoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D30") "5.05"
Set-TextValue $ws.Range("E30") "  +7.31%  "

Set-TextValue $ws.Range("B31") "ImmutableX"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D31") "1.18"
Set-TextValue $ws.Range("E31") "  +1.37%  "

# Row 39 <-> Row 40: HuobiToken and THORChain swapped ranking order, with
# refreshed price/volume figures.
Set-TextValue $ws.Range("B39") "THORChain"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D39") "5.15"
Set-TextValue $ws.Range("E39") "  +9.98%  "

Set-TextValue $ws.Range("B40") "HuobiToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D40") "3.08"
Set-TextValue $ws.Range("E40") "  +8.52%  "
